$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the existing "总计" sheet to "2022-Q1" and duplicate it
#    (before touching its data) so the duplicate inherits all of its
#    formatting (sheetPr / pageMargins / styles) for the new "总计"
#    sheet that will be re-created further below.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

$q1.Copy($null, $q1)
$total = $wb.Worksheets.Item("2022-Q1 (2)")
$total.Name = "总计"

# ---------------------------------------------------------------------
# 2. Rewrite the "2022-Q1" sheet with the new fund-holding table
#    (columns A-H, same layout used by the 2021-Q3 / 2021-Q4 sheets).
# ---------------------------------------------------------------------

# Extend the header row's styling (s="2") from column D into E:H first,
# since those cells do not exist yet in the old sheet.
$q1.Range("D1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

$q1.Cells.Item(2,1).Value = 0
$q1.Cells.Item(2,2).Value = "'001128"
$q1.Cells.Item(2,3).Value = "宝盈新兴产业灵活配置混合"
$q1.Cells.Item(2,4).Value = "'19.94"
$q1.Cells.Item(2,5).Value = "'86.76"
$q1.Cells.Item(2,6).Value = "'3.50"
$q1.Cells.Item(2,7).Value = "'0.6979"
$q1.Cells.Item(2,8).Value = 8

$q1.Cells.Item(3,1).Value = 1
$q1.Cells.Item(3,2).Value = "'420005"
$q1.Cells.Item(3,3).Value = "天弘周期策略混合"
$q1.Cells.Item(3,4).Value = "'5.25"
$q1.Cells.Item(3,5).Value = "'89.31"
$q1.Cells.Item(3,6).Value = "'4.96"
$q1.Cells.Item(3,7).Value = "'0.2604"
$q1.Cells.Item(3,8).Value = 9

$q1.Cells.Item(4,1).Value = 2
$q1.Cells.Item(4,2).Value = "'007202"
$q1.Cells.Item(4,3).Value = "天弘优质成长企业精选混合"
$q1.Cells.Item(4,4).Value = "'4.81"
$q1.Cells.Item(4,5).Value = "'92.52"
$q1.Cells.Item(4,6).Value = "'4.10"
$q1.Cells.Item(4,7).Value = "'0.1972"
$q1.Cells.Item(4,8).Value = 10

# ---------------------------------------------------------------------
# 3. Rewrite the "总计" sheet with the updated summary table, adding
#    the new 2022-Q1 row on top and re-numbering the A column index.
# ---------------------------------------------------------------------

# Extend column A's index styling (s="2") down into the newly added
# 5th row, which does not exist yet in the copied sheet.
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$total.Cells.Item(1,2).Value = "日期"
$total.Cells.Item(1,3).Value = "持有数量(只)"
$total.Cells.Item(1,4).Value = "持有市值(亿元)"

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 3
$total.Cells.Item(2,4).Value = 1.16

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2021-Q4"
$total.Cells.Item(3,3).Value = 2
$total.Cells.Item(3,4).Value = 0.48

$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2021-Q3"
$total.Cells.Item(4,3).Value = 2
$total.Cells.Item(4,4).Value = 0.43

$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(5,2).Value = "2021-Q1"
$total.Cells.Item(5,3).Value = 1
$total.Cells.Item(5,4).Value = 0.02

# ---------------------------------------------------------------------
# 4. Restore the originally active sheet.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q1").Activate()
